# Commit: "Added More Certain Missing Elements" / "Good Job Luke -.-"
#
# Functional edits captured by the diff:
#  1. Cell D3 on Sheet1 changes from 14.8104 to 10.2342 (the shared formula
#     in I3 recalculates automatically as a result: 11.43178 -> 10.516540000000001).
#  2. The active selection on Sheet1 moves from A1:I5 (active cell I1) to
#     the single cell J11.
#
# (The remaining hunks in the source diff - fileVersion/calcId bumps, the
# x15 absPath, extra theme fonts, chart c:layout/extLst/uniqueId churn, and
# the a16:creationId blocks on the drawing shapes - are artifacts written by
# a newer Excel build when it resaves the package and aren't reachable via
# the Excel object model, so they aren't reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the data value; dependent formulas recalc automatically.
$ws.Range("D3").Value = 10.2342

# 2. Move the selection to J11.
$ws.Range("J11").Select()
